$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''245.76'
$ws.Range('E2').Value = '''-0.36%'
$ws.Range('D3').Value = '''30.02'
$ws.Range('E3').Value = '''-2.56%'
$ws.Range('D4').Value = '''5.156'
$ws.Range('E4').Value = '''-0.35%'
$ws.Range('E5').Value = '''0.46%'
$ws.Range('D6').Value = '''6.668'
$ws.Range('E6').Value = '''0.99%'
$ws.Range('D7').Value = '''3.292'
$ws.Range('E7').Value = '''7.19%'
$ws.Range('D8').Value = '''0.8503'
$ws.Range('E8').Value = '''-0.55%'
$ws.Range('D9').Value = '''0.8576'
$ws.Range('E9').Value = '''-2.55%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '''0.1385'
$ws.Range('E10').Value = '''1.33%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '''0.07077'
$ws.Range('E11').Value = '''-0.05%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '''0.03232'
$ws.Range('E12').Value = '''12.94%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '''0.09363'
$ws.Range('E13').Value = '''-0.33%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '''0.001535'
$ws.Range('E14').Value = '''1.09%'
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').Value = '''0.0005970'
$ws.Range('E15').Value = '''-94.21%'
$ws.Range('D16').Value = '''0.006016'
$ws.Range('E16').Value = '''-1.33%'
$ws.Range('D17').Value = '''3.548'
$ws.Range('E17').Value = '''1.68%'
$ws.Range('D18').Value = '''2.182'
$ws.Range('E18').Value = '''0.29%'
$ws.Range('E19').Value = '''-2.38%'
$ws.Range('D20').Value = '''0.03418'
$ws.Range('E20').Value = '''5.40%'
$ws.Range('D21').Value = '''0.1323'
$ws.Range('E21').Value = '''1.82%'
$ws.Range('D22').Value = '''3.491'
$ws.Range('E22').Value = '''-1.04%'
$ws.Range('B23').Value = 'CoinExToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D23').Value = '''0.04134'
$ws.Range('E23').Value = '''-0.50%'
$ws.Range('B24').Value = 'ZBToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D24').Value = '''0.1410'
$ws.Range('E24').Value = '''2.22%'
$ws.Range('D25').Value = '''0.001224'
$ws.Range('E25').Value = '''0.65%'
$ws.Range('D26').Value = '''0.004157'
$ws.Range('E26').Value = '''-7.54%'
$ws.Range('E28').Value = '''4.83%'
$ws.Range('D40').Value = '''0.03746'
$ws.Range('E40').Value = '''-0.80%'
$ws.Range('D41').Value = '''0.1070'
$ws.Range('E41').Value = '''-0.07%'
$ws.Range('B42').Value = 'KickToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D42').Value = '''0.003582'
$ws.Range('E42').Value = '''-36.91%'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').Value = '''0.002460'
$ws.Range('E43').Value = '''-5.34%'
$ws.Range('D44').Value = '''0.01070'
$ws.Range('E44').Value = '''6.88%'
$ws.Range('D45').Value = '''0.00005484'
$ws.Range('E45').Value = '''8.10%'
$ws.Range('E46').Value = '''0.08%'
$ws.Range('D47').Value = '''0.07100'
$ws.Range('E47').Value = '''-11.18%'
$ws.Range('D48').Value = '''0.002474'
$ws.Range('E48').Value = '''-10.48%'
$ws.Range('D49').Value = '''0.00002100'
$ws.Range('E49').Value = '''0.08%'
$ws.Range('D50').Value = '''0.0002000'
$ws.Range('E50').Value = '''0.08%'
